$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) figures increased by 1 for two events.
# These totals live in column F on both the "展览" sheet and the
# "全部类型" aggregate sheet, which mirrors the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 357
    $ws.Range("F5").Value = 294
}
